$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "36.524.76"
$ws.Range("E2").Value = "  +2.90%  "

# Row 3
$ws.Range("D3").Value = "2.068.14"
$ws.Range("E3").Value = "  +9.68%  "

# Row 4
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.39%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.663"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.78%  "

# Row 7
$ws.Range("E7").Value = "  +0.11%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.99"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.05%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "60.64"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.28%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.365"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.99%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0718"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.33%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0986"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.26%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.44"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.91%  "

# Row 14
$ws.Range("D14").Value = "2.372.95"
$ws.Range("E14").Value = "  +9.78%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.815"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.10%  "

# Row 16
$ws.Range("D16").Value = "2.062.85"
$ws.Range("E16").Value = "  +9.31%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.90"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.30%  "

# Row 18
$ws.Range("D18").Value = "36.504.67"
$ws.Range("E18").Value = "  +2.84%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.77%  "

# Row 20
$ws.Range("E20").Value = "  -1.97%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "237.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.29%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.60"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.05%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.89"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.26%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.11%  "

# Row 25
$ws.Range("E25").Value = "  -7.35%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.02"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.07%  "

# Row 27
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.00%  "

# Row 28
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +9.71%  "

# Row 29
$ws.Range("E29").Value = "  -8.51%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.121"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.82%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "21.60"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +51.04%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.19%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0580"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.12%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0897"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +19.93%  "

# Row 35
$ws.Range("E35").Value = "  +0.12%  "

# Row 36
$ws.Range("E36").Value = "  +0.54%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.24"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +15.74%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.99"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.50%  "

# Row 39
$ws.Range("E39").Value = "  +2.44%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.32"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -10.52%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.13"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.68%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "96.62"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.07%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0214"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.68%  "

# Row 44
$ws.Range("E44").Value = "  +15.79%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.94"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.05%  "

# Row 46
$ws.Range("D46").Value = "1.315.97"
$ws.Range("E46").Value = "  +0.56%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0817"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.80%  "

# Row 49
$ws.Range("D49").Value = "2.250.61"
$ws.Range("E49").Value = "  +9.12%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.16%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.83"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +15.13%  "
